$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): a new "Odd_CS_3-3_HT" header is inserted at AW1,
# shifting the previous AW1:BC1 headers one column to the right (BD1
# "Odd_CS_4-4_HT" is unaffected). Apply the resulting header values directly.
$ws.Range("AW1").Value = "Odd_CS_3-3_HT"
$ws.Range("AX1").Value = "Odd_CS_0-1_HT"
$ws.Range("AY1").Value = "Odd_CS_0-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_1-2_HT"
$ws.Range("BA1").Value = "Odd_CS_0-3_HT"
$ws.Range("BB1").Value = "Odd_CS_1-3_HT"
$ws.Range("BC1").Value = "Odd_CS_2-3_HT"

# --- Data row 2: fully updated match info/odds.
$ws.Range("A2").Value = "ttLF3hdB"
$ws.Range("B2").Value = "30/10/2024"
$ws.Range("C2").Value = "07:00"
$ws.Range("D2").Value = "JAPAN - J1 LEAGUE"
$ws.Range("E2").Value = "Yokohama F. Marinos"
$ws.Range("F2").Value = "Urawa Reds"
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 2.9
$ws.Range("J2").Value = 2.88
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.35
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 17
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 13
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 29
$ws.Range("AL2").Value = 21
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.75
$ws.Range("AO2").Value = 12
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 41
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 5
$ws.Range("AY2").Value = 15
$ws.Range("AZ2").Value = 19
$ws.Range("BA2").Value = 41
$ws.Range("BB2").Value = 51
$ws.Range("BC2").Value = 101
$ws.Range("BD2").Value = 81
